$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.907.62"
$ws.Range("D3").Value = "1.551.24"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("E4").Value = "  -0.54%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.54"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("E7").Value = "  -0.54%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.00"
$ws.Range("E8").Value = "  +1.77%  "
$ws.Range("E9").Value = "  -0.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0594"
$ws.Range("E10").Value = "  +0.89%  "
$ws.Range("E11").Value = "  -0.63%  "
$ws.Range("D12").Value = "1.772.64"
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("D13").Value = "1.549.84"
$ws.Range("E13").Value = "  -0.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.74"
$ws.Range("E14").Value = "  +0.70%  "
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").Value = "26.905.66"
$ws.Range("E16").Value = "  -0.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.58"
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("D18").Value = "0.0₃0708"
$ws.Range("E18").Value = "  +2.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "217.05"
$ws.Range("E19").Value = "  +0.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.28"
$ws.Range("E20").Value = "  +0.20%  "
$ws.Range("E21").Value = "  -0.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.07"
$ws.Range("E22").Value = "  +1.00%  "
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("E24").Value = "  -1.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.41"
$ws.Range("E25").Value = "  +0.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.63"
$ws.Range("E26").Value = "  -0.38%  "
$ws.Range("E27").Value = "  +0.49%  "
$ws.Range("E28").Value = "  +0.67%  "
$ws.Range("E29").Value = "  -0.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0469"
$ws.Range("E30").Value = "  +1.26%  "
$ws.Range("E31").Value = "  -1.22%  "
$ws.Range("E32").Value = "  -0.37%  "
$ws.Range("E33").Value = "  +3.71%  "
$ws.Range("D34").Value = "1.414.47"
$ws.Range("E34").Value = "  +0.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.60"
$ws.Range("E35").Value = "  +2.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.985"
$ws.Range("E36").Value = "  +3.50%  "
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0165"
$ws.Range("E38").Value = "  +0.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.528"
$ws.Range("E39").Value = "  +0.88%  "
$ws.Range("E40").Value = "  -0.50%  "
$ws.Range("E41").Value = "  -0.56%  "
$ws.Range("E42").Value = "  +2.57%  "
$ws.Range("E43").Value = "  +1.67%  "
$ws.Range("E44").Value = "  +0.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.39"
$ws.Range("E45").Value = "  +0.88%  "
$ws.Range("E46").Value = "  -0.73%  "
$ws.Range("D47").Value = "1.686.07"
$ws.Range("E47").Value = "  -0.27%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.15"
$ws.Range("E48").Value = "  +1.01%  "
$ws.Range("E49").Value = "  +1.68%  "
$ws.Range("E50").Value = "  +2.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0958"
$ws.Range("E51").Value = "  +0.13%  "
